$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Column G holds the "Type" label for each measurement (e.g.
# "swab_SARSCoV2_E_Ct", "sputum_SARSCoV2_RdRp_VL", ...). Strip the
# trailing "_Ct" / "_VL" suffix from every data row so the label just
# reflects target + specimen (e.g. "swab_SARSCoV2_RdRp").
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Text
    if ($val -ne $null) {
        $newVal = $val -replace '_Ct$', '' -replace '_VL$', ''
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# Update the view state to match: scrolled so row 105 is at top, selection G127
$ws.Range("G127").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 105
$win.ScrollColumn = 1
